# Weekly update: a new daily record is inserted at the top of the
# "Vega Modelo de Temuco - Pepino dulce" data block (row 88), pushing the
# existing historical rows (88-174) down by one (now 89-175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 88; everything below shifts down one row.
$ws.Rows("88:88").Insert()

# Populate the newly inserted row with the new day's reading.
$ws.Range("A88").Value = 10
$ws.Range("B88").Value = "Vega Modelo de Temuco"
$ws.Range("C88").Value = "La Araucanía"
$ws.Range("D88").Value = 44587
$ws.Range("E88").Value = 9
$ws.Range("F88").Value = 100112043
$ws.Range("G88").Value = "Pepino dulce"
$ws.Range("H88").Value = "Cultivar IV Región"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 35
$ws.Range("K88").Value = 22000
$ws.Range("L88").Value = 22000
$ws.Range("M88").Value = 22000
$ws.Range("N88").Value = "`$/bandeja 18 kilos"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 1222
$ws.Range("Q88").Value = 18
$ws.Range("R88").Value = "Hortaliza"
